$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original values for columns D, M, N, O, P, R, S across rows 2-31
# before overwriting anything, since the update is a permutation of rows.
$snapD = @{}
$snapM = @{}
$snapN = @{}
$snapO = @{}
$snapP = @{}
$snapR = @{}
$snapS = @{}

for ($r = 2; $r -le 31; $r++) {
    $snapD[$r] = $ws.Cells.Item($r, 4).Value()
    $snapM[$r] = $ws.Cells.Item($r, 13).Value()
    $snapN[$r] = $ws.Cells.Item($r, 14).Value()
    $snapO[$r] = $ws.Cells.Item($r, 15).Value()
    $snapP[$r] = $ws.Cells.Item($r, 16).Value()
    $snapR[$r] = $ws.Cells.Item($r, 18).Value()
    $snapS[$r] = $ws.Cells.Item($r, 19).Value()
}

# Mapping: destination row -> source row (values originally in source row now belong to destination row)
$rowMap = @{}
$rowMap[2] = 5
$rowMap[3] = 6
$rowMap[4] = 25
$rowMap[5] = 11
$rowMap[6] = 9
$rowMap[7] = 10
$rowMap[8] = 15
$rowMap[9] = 23
$rowMap[10] = 29
$rowMap[11] = 17
$rowMap[12] = 14
$rowMap[13] = 2
$rowMap[14] = 4
$rowMap[15] = 18
$rowMap[16] = 21
$rowMap[17] = 12
$rowMap[18] = 8
$rowMap[19] = 24
$rowMap[20] = 13
$rowMap[21] = 26
$rowMap[22] = 7
$rowMap[23] = 3
$rowMap[24] = 22
$rowMap[25] = 31
$rowMap[26] = 30
$rowMap[27] = 20
$rowMap[28] = 19
$rowMap[29] = 16
$rowMap[30] = 28
$rowMap[31] = 27

foreach ($dst in $rowMap.Keys) {
    $src = $rowMap[$dst]
    $ws.Cells.Item($dst, 4).Value = $snapD[$src]
    $ws.Cells.Item($dst, 13).Value = $snapM[$src]
    $ws.Cells.Item($dst, 14).Value = $snapN[$src]
    $ws.Cells.Item($dst, 15).Value = $snapO[$src]
    $ws.Cells.Item($dst, 16).Value = $snapP[$src]
    $ws.Cells.Item($dst, 18).Value = $snapR[$src]
    $ws.Cells.Item($dst, 19).Value = $snapS[$src]
}

Write-Host "Done applying row permutation."